$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2: DummyName -> Abangan Norte
$ws.Range("A2").Value = "Abangan Norte"
Set-TextValue "B2" "14.7661"
Set-TextValue "C2" "120.9342"
Set-TextValue "D2" "11417"
Set-TextValue "E2" "10080"
Set-TextValue "F2" "6431"

# Row 3: Abangan -> Abangan Sur
$ws.Range("A3").Value = "Abangan Sur"
Set-TextValue "B3" "14.7653"
Set-TextValue "C3" "120.9437"
Set-TextValue "D3" "10595"
Set-TextValue "E3" "9750"
Set-TextValue "F3" "6525"

# Row 4: Hotdog -> Ibayo
$ws.Range("A4").Value = "Ibayo"
Set-TextValue "B4" "14.7535"
Set-TextValue "C4" "120.9533"
Set-TextValue "D4" "8310"
Set-TextValue "E4" "5000"
Set-TextValue "F4" "7186"
